$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values, entered in the order that reproduces the shared-string table order ---
$ws.Range("A1").Value = "Test Case Name"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Token"
$ws.Range("A2").Value = "Login_to_Member_Portal"
$ws.Range("A3").Value = "Logout_of_Member_Portal"
$ws.Range("F1").Value = "ResetHeader"
$ws.Range("F4").Value = "Member Security Questions"
$ws.Range("B4").Value = "testuserd3432"
$ws.Range("C4").Value = "tDuestest%3432"
$ws.Range("E1").Value = "NewPassword"
$ws.Range("E4").Value = "Test123"
$ws.Range("A4").Value = "First_Time_Reset_Password"
$ws.Range("H1").Value = "Q1Ans"
$ws.Range("J1").Value = "Q2Ans"
$ws.Range("I1").Value = "Q2Val"
$ws.Range("G1").Value = "Q1Val"
$ws.Range("H4").Value = "Brown"
$ws.Range("J4").Value = "Texas"
$ws.Range("G4").Value = "'436"
$ws.Range("I4").Value = "'439"

$ws.Range("B2").Value = "test123"
$ws.Range("C2").Value = "'123"
$ws.Range("D2").Value = "XT13"
$ws.Range("B3").Value = "test123"
$ws.Range("C3").Value = "'123"
$ws.Range("D3").Value = "XT13"
$ws.Range("D4").Value = "XT13"

# --- Column widths (best-effort visual match) ---
$ws.Columns("A:A").ColumnWidth = 25.14
$ws.Columns("B:B").ColumnWidth = 10.71
$ws.Columns("C:C").ColumnWidth = 12.28
$ws.Columns("D:D").ColumnWidth = 12.71
$ws.Columns("E:E").ColumnWidth = 15.42
$ws.Columns("F:F").ColumnWidth = 22.14

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("J1").Select()
